# "Server Tests Setup" remedies workbook - add a dedicated resource/setting
# remedy for the ResumingAWorkflowWithAnInvalidUserReturnsAuthenticationError
# spec (row 8), matching the other per-spec "copy secure.config" remedies but
# pointing at its own config sub-directory.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C8 previously shared the generic "Warewolf.Security.Specs_Setup" secure.config
# remedy text (shared string 14); give this row its own specific instructions,
# which appends a brand-new shared string to the table.
$ws.Range("C8").Value = 'Copy secure.config (found at the "ResumingAWorkflowWithAnInvalidUserReturnsAuthenticationError config" directory) file in the Directory "%Programdata%\Warewolf\Server Settings"'

# The extra text now wraps onto a third line at the current column width, so
# the row grows from two lines (28.8pt) to three (43.2pt).
$ws.Rows.Item(8).RowHeight = 43.2

# Reflect the saved file's scroll position/selection (view was left scrolled
# to column B with C8 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("C8").Select()
